# Lots of bugfixes & improvements
$wb = $excel.ActiveWorkbook

# --- "Item Levels" sheet: tweak the level-scaling constants ---------------
$wsItem = $wb.Worksheets.Item("Item Levels")
$wsItem.Range("H2").Value = 1.5
$wsItem.Range("H3").Value = 1

# Leave this sheet's selection at H4 (it is no longer the active tab).
$wsItem.Activate() | Out-Null
$wsItem.Range("H4").Select() | Out-Null

# --- "Map Size" sheet becomes the active / selected tab -------------------
$wsMap = $wb.Worksheets.Item("Map Size")
$wsMap.Activate() | Out-Null
$wsMap.Range("F22").Select() | Out-Null
